# OSB Reminder Form - basic pdf generation template tweaks
$d = $word.ActiveDocument

# 1. Title paragraph: add space-before (12pt = 240 twips) to the
#    "Assistance League of the Eastside" paragraph.
$p1 = $d.Paragraphs(1)
$p1.SpaceBefore = 12

# 2. Reposition the anchored logo picture further down
#    (50165 EMU -> 173990 EMU, i.e. 3.95pt -> 13.7pt).
$shape = $d.Shapes(1)
$shape.Top = 13.7

# 3. Merge the "Assistance League " / "of the Eastside" runs (which were
#    split by a stray _GoBack bookmark) into a single run.
$d.Content.Find.Execute("Assistance League of the Eastside", $false, $false, `
    $false, $false, $false, $true, 1, $false, `
    "Assistance League of the Eastside", 2) | Out-Null

# 4. "Shopping Event Information" heading: bump space-before 12pt -> 24pt.
$p3 = $d.Paragraphs(3)
$p3.SpaceBefore = 24

# 5. Body copy + bullet list (paragraphs 4-11): shrink 14pt -> 13pt and
#    stamp matching complex-script size.
for ($i = 4; $i -le 11; $i++) {
    $rng = $d.Paragraphs($i).Range
    $rng.Font.Size = 13
    $rng.Font.SizeBi = 13
}

# 6. Last bullet ("All purchases are final...") space-after 12pt -> 18pt.
$p11 = $d.Paragraphs(11)
$p11.SpaceAfter = 18

# 7. Split the "clothing and shoes and may not purchase..." run and move
#    the _GoBack bookmark in between, right after "shoes".
$p8 = $d.Paragraphs(8)
$splitPoint = $p8.Range.Start + 70
$bookmarkRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

# 8. Event-Information table header cell: top padding 4.3pt -> 7.2pt.
$table = $d.Tables(1)
$table.Cell(1, 1).TopPadding = 7.2

# 9. Blank spacer rows inside the table (7pt -> 9pt), skipping the
#    centered one in the very first spacer row which stays at 7pt.
#    (Iterate via table cells rather than Paragraphs - COM Range bounds
#    on a wholly-empty paragraph don't resolve reliably in this host.)
for ($r = 1; $r -le $table.Rows.Count; $r++) {
    $cell = $table.Cell($r, 1)
    if ($cell.Range.Font.Size -eq 7 -and $cell.Range.ParagraphFormat.Alignment -ne 1) {
        $cell.Range.Font.Size = 9
    }
}
